# Juno: check in to OLPRODLOC.
# Translate the sales report workbook (English -> Italian).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Report delle vendite"

# Header row translations.
$ws.Range("A1").Value = "Anno-Trimestre"
$ws.Range("F1").Value = "sud-orientale"

# Year-Quarter column translations (Q -> T for "Trimestre").
$ws.Range("A2").Value = "2022-T1"
$ws.Range("A3").Value = "2022-T2"
$ws.Range("A4").Value = "2022-T3"
$ws.Range("A5").Value = "2022-T4"
$ws.Range("A6").Value = "2023-T1"
$ws.Range("A7").Value = "2023-T2"
$ws.Range("A8").Value = "2023-T3"
$ws.Range("A9").Value = "2023-T4"
